# Removed cc to support email
# Append " *" to each header label in row 1 (Fund, Title, Gross, Carry, Date)
# and move the active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Fund *"
$ws.Range("B1").Value = "Title *"
$ws.Range("C1").Value = "Gross *"
$ws.Range("D1").Value = "Carry *"
$ws.Range("E1").Value = "Date *"

$ws.Range("A2").Select()
